# Adds a "Time" column and an "Average UPH" summary row to the
# PUTWALL PICKING and REGULAR PICK sheets, re-orders the data rows by
# descending UPH, and highlights the header / summary rows with a light
# blue fill.

$wb = $excel.ActiveWorkbook

$lightBlue = 15128749  # RGB(173, 216, 230) -> ADD8E6

function Fill-Sheet {
    param(
        $ws,
        [string]$qtyHeader,
        $rows,
        [double]$avgUph
    )

    # Header row
    $ws.Cells.Item(1, 1).Value = "UserID"
    $ws.Cells.Item(1, 2).Value = $qtyHeader
    $ws.Cells.Item(1, 3).Value = "Time"
    $ws.Cells.Item(1, 4).Value = "UPH"

    # Data rows (re-ordered, with the new Time column and recalculated UPH)
    $r = 2
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $r = $r + 1
    }

    # Summary row
    $ws.Cells.Item($r, 1).Value = "Average UPH"
    $ws.Cells.Item($r, 4).Value = $avgUph

    $lastCol = "D"
    $ws.Range("A1:" + $lastCol + "1").Interior.Color = $lightBlue
    $ws.Range("A" + $r + ":" + $lastCol + $r).Interior.Color = $lightBlue
}

# ---------------------------------------------------------------------------
# PUTWALL PICKING
# ---------------------------------------------------------------------------
$wsPutwall = $wb.Worksheets.Item("PUTWALL PICKING")

$putwallRows = @(
    @("DIAN4065.ENTRIALGO", 197, 44, 268.64),
    @("ANASTASIIA.MAKHTOUT", 276, 91, 181.98),
    @("ABHI4088.ABHISHEK", 286, 104, 165),
    @("KADE3054.ZONGO", 169, 67, 151.34),
    @("BOHD0676.KUSHLIAK", 173, 76, 136.58),
    @("LOWRHY-OTIENO.JAOKO", 120, 53, 135.85),
    @("TANI2739.HOSSAINISLA", 378, 167, 135.81),
    @("STAN9294.BAUER", 94, 46, 122.61),
    @("THIE6554.DIALLO", 351, 182, 115.71),
    @("WILDINE.JEUNE", 141, 88, 96.14),
    @("RAVI4279.THAKUR", 29, 10, 0),
    @("SEPIDEH.AZARIHASHJIN", 24, 6, 0),
    @("MDSAIFUL.ISLAM", 66, 19, 0),
    @("RARG046N.YEBOAH", 3, 0, 0),
    @("NESR2403.ATTALAH", 49, 18, 0),
    @("AGNE8120.CARUTH", 12, 19, 0),
    @("LOANA.MBONGO", 1, 0, 0),
    @("HARJ4282.SINGH", 51, 17, 0),
    @("DEVI789.SINGH", 46, 15, 0),
    @("YATI0689.YATIN", 39, 9, 0)
)

Fill-Sheet $wsPutwall "PutwallPickingQuantity" $putwallRows 150.97

# ---------------------------------------------------------------------------
# REGULAR PICK
# ---------------------------------------------------------------------------
$wsRegular = $wb.Worksheets.Item("REGULAR PICK")

$regularRows = @(
    @("DIAN4065.ENTRIALGO", 179, 38, 282.63),
    @("BOHD0676.KUSHLIAK", 148, 110, 80.73),
    @("WILDINE.JEUNE", 22, 43, 30.7),
    @("SEPIDEH.AZARIHASHJIN", 14, 34, 24.71),
    @("AGNE8120.CARUTH", 8, 3, 0),
    @("NESR2403.ATTALAH", 5, 0, 0),
    @("ZAHIDGUL.MINHAS", 17, 19, 0),
    @("THIE6554.DIALLO", 17, 22, 0),
    @("TANI2739.HOSSAINISLA", 2, 0, 0),
    @("SURESH.DHAWAN", 2, 0, 0),
    @("STAN9294.BAUER", 18, 20, 0),
    @("RAVI4279.THAKUR", 2, 0, 0),
    @("RARG046N.YEBOAH", 20, 1, 0),
    @("MARI882N.ABDELKADER", 17, 15, 0),
    @("ANASTASIIA.MAKHTOUT", 17, 4, 0),
    @("LOWRHY-OTIENO.JAOKO", 35, 23, 0),
    @("LOANA.MBONGO", 1, 0, 0),
    @("KHINEHAYMAR.THAUNG", 1, 0, 0),
    @("KADE3054.ZONGO", 1, 0, 0),
    @("JEEW9554.SITUMUDALIG", 5, 3, 0),
    @("HARJ4282.SINGH", 2, 0, 0),
    @("DEVI789.SINGH", 1, 0, 0),
    @("ARJUNBHAI.PATEL", 27, 7, 0),
    @("ZAKI0190.PHILLIPHORS", 9, 12, 0)
)

Fill-Sheet $wsRegular "RegularPickQuantity" $regularRows 104.69
